$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 — ZERNAKH ABDELLAH replaces KJKKOPPP JKLML
$ws.Range("A2").Value = "ZERNAKH ABDELLAH"
$ws.Range("B2").Value = "IB19558"
# Long digit-only account number: prefix with an apostrophe so Excel keeps
# it as text instead of coercing it to a (precision-losing) number.
$ws.Range("C2").Value = "'145101211406073828000084"
$ws.Range("D2").Value = "MARRAKECH BENI MELLAL"
$ws.Range("E2").Value = "BP"
$ws.Range("F2").Value = "Point de vente"
$ws.Range("G2").Value = "052/FKIH BEN SALEH"
$ws.Range("H2").Value = "mensuelle"
$ws.Range("I2").Value = 11000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 11000

# Row 3 — NOUBAIL MOUNTASSIR replaces VAVA FAFA
$ws.Range("A3").Value = "NOUBAIL MOUNTASSIR"
$ws.Range("B3").Value = "Q251990"
$ws.Range("C3").Value = "'007400000313200019604463"
$ws.Range("D3").Value = "KHOURIBGA ZELLAKA"
$ws.Range("E3").Value = "AWB"
$ws.Range("F3").Value = "Direction régionale"
$ws.Range("G3").Value = "905/TADLA OUARDIGHA ZAYANE"
$ws.Range("H3").Value = "mensuelle"
$ws.Range("I3").Value = 6750
$ws.Range("J3").Value = 675
$ws.Range("K3").Value = 6075

# Row 4 — NOUBAIL MOHAMMED replaces CHARIJI ABDELLAH
$ws.Range("A4").Value = "NOUBAIL MOHAMMED"
$ws.Range("B4").Value = "IR801997"
$ws.Range("C4").Value = "'007400000313200019604463"
$ws.Range("D4").Value = "KHOURIBGA ZELLAKA"
$ws.Range("E4").Value = "AWB"
$ws.Range("F4").Value = "Direction régionale"
$ws.Range("G4").Value = "905/TADLA OUARDIGHA ZAYANE"
$ws.Range("H4").Value = "mensuelle"
$ws.Range("I4").Value = 6750
$ws.Range("J4").Value = 675
$ws.Range("K4").Value = 6075

# Row 5 used to hold "JALAL MED"; it now becomes the blank totals row that
# used to be row 6 (rows 2-4 are summed instead of rows 2-5).
$ws.Range("A5").Value = " "
$ws.Range("B5").Value = " "
$ws.Range("C5").Value = " "
$ws.Range("D5").Value = " "
$ws.Range("E5").Value = " "
$ws.Range("F5").Value = " "
$ws.Range("G5").Value = " "
$ws.Range("H5").Value = " "
$ws.Range("I5").Value = 24500
$ws.Range("J5").Value = 1350
$ws.Range("K5").Value = 23150

# The old row 6 (previous totals row) is removed entirely, shrinking the
# sheet's used range down to A1:K5.
$ws.Range("A6:K6").EntireRow.Delete()
